$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.945
$ws.Range("A21").Value = -19.953
$ws.Range("A23").Value = -20.131
$ws.Range("A25").Value = -21.742
$ws.Range("A53").Value = -21.928
$ws.Range("A57").Value = -22.253
$ws.Range("A59").Value = -22.461
$ws.Range("A69").Value = -21.519
$ws.Range("A79").Value = -21.073
$ws.Range("A83").Value = -21.902
$ws.Range("A93").Value = -21.508
